$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new data row (row 29), copying formatting (incl. the date
# number format) from the last existing data row (row 28) so no new
# style entries are minted.
$ws.Range("A28").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = 43919

$ws.Range("B29").Value = 1299
$ws.Range("C29").Value = 96
$ws.Range("D29").Value = 28
$ws.Range("E29").Value = 68
$ws.Range("F29").Value = 0

# Update the active selection to match the newly added last cell
$ws.Range("F29").Select()
